# Update crafting-leve profit metrics (currentAveragePrice* / LevePrice* / LeveProfit*)
# for the rows whose market data changed, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 1199712.5
$ws.Range("I86").Value = 1797767.5
$ws.Range("J86").Value = 3602.5557
$ws.Range("K86").Value = 1797767.5
$ws.Range("L86").Value = 3602.5557
$ws.Range("M86").Value = -1796644.5
$ws.Range("N86").Value = -5848.5557
# Row 89
$ws.Range("H89").Value = 1199712.5
$ws.Range("I89").Value = 1797767.5
$ws.Range("J89").Value = 3602.5557
$ws.Range("K89").Value = 8988837.5
$ws.Range("L89").Value = 18012.7785
$ws.Range("M89").Value = -8983221.5
$ws.Range("N89").Value = -29244.7785
# Row 96
$ws.Range("H96").Value = 5066.9414
$ws.Range("J96").Value = 10392.125
$ws.Range("L96").Value = 31176.375
$ws.Range("N96").Value = -33922.375
# Row 113
$ws.Range("H113").Value = 4750
$ws.Range("J113").Value = 4750
$ws.Range("L113").Value = 4750
$ws.Range("N113").Value = -11258
# Row 121
$ws.Range("H121").Value = 3959.875
$ws.Range("J121").Value = 3959.875
$ws.Range("L121").Value = 11879.625
$ws.Range("N121").Value = -15373.625
# Row 132
$ws.Range("H132").Value = 12119.8545
$ws.Range("I132").Value = 5295.9375
$ws.Range("J132").Value = 14919.41
$ws.Range("K132").Value = 15887.8125
$ws.Range("L132").Value = 44758.23
$ws.Range("M132").Value = -13357.8125
$ws.Range("N132").Value = -49818.23
# Row 137
$ws.Range("H137").Value = 6805952
$ws.Range("I137").Value = 1014.2692
$ws.Range("J137").Value = 14498490
$ws.Range("K137").Value = 3042.8076
$ws.Range("L137").Value = 43495470
$ws.Range("M137").Value = -492.8076000000001
$ws.Range("N137").Value = -43500570
# Row 138
$ws.Range("H138").Value = 3974.4375
$ws.Range("J138").Value = 4307.7
$ws.Range("L138").Value = 12923.1
$ws.Range("N138").Value = -23203.1

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4241.0845
$ws.Range("I32").Value = 2262.7856
$ws.Range("K32").Value = 2262.7856
$ws.Range("M32").Value = -1975.7856
# Row 102
$ws.Range("H102").Value = 298969.97
$ws.Range("I102").Value = 508302.38
$ws.Range("K102").Value = 508302.38
$ws.Range("M102").Value = -506680.38
# Row 122
$ws.Range("H122").Value = 3960.75
$ws.Range("I122").Value = 2336.3076
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 7008.9228
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -4558.9228
$ws.Range("N122").Value = -37900
# Row 132
$ws.Range("H132").Value = 21125.438
$ws.Range("I132").Value = 19548.838
$ws.Range("J132").Value = 70000
$ws.Range("K132").Value = 58646.514
$ws.Range("L132").Value = 210000
$ws.Range("M132").Value = -56116.514
$ws.Range("N132").Value = -215060

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 58
$ws.Range("H58").Value = 37029.6
$ws.Range("I58").Value = 37173
$ws.Range("J58").Value = 36814.5
$ws.Range("K58").Value = 37173
$ws.Range("L58").Value = 36814.5
$ws.Range("M58").Value = -36879
$ws.Range("N58").Value = -37402.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3131.9048
$ws.Range("I31").Value = 998.73334
$ws.Range("K31").Value = 998.73334
$ws.Range("M31").Value = -703.73334
# Row 34
$ws.Range("H34").Value = 3131.9048
$ws.Range("I34").Value = 998.73334
$ws.Range("K34").Value = 998.73334
$ws.Range("M34").Value = -796.73334
# Row 132
$ws.Range("H132").Value = 22239588
$ws.Range("I132").Value = 27797320
$ws.Range("J132").Value = 8665
$ws.Range("K132").Value = 83391960
$ws.Range("L132").Value = 25995
$ws.Range("M132").Value = -83389430
$ws.Range("N132").Value = -31055
# Row 134
$ws.Range("H134").Value = 3093.3096
$ws.Range("I134").Value = 3103
$ws.Range("J134").Value = 2899.5
$ws.Range("K134").Value = 9309
$ws.Range("L134").Value = 8698.5
$ws.Range("M134").Value = -6774
$ws.Range("N134").Value = -13768.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 2715.6667
$ws.Range("J25").Value = 3998.5
$ws.Range("L25").Value = 11995.5
$ws.Range("N25").Value = -12333.5
# Row 30
$ws.Range("H30").Value = 2715.6667
$ws.Range("J30").Value = 3998.5
$ws.Range("L30").Value = 11995.5
$ws.Range("N30").Value = -12199.5
# Row 113
$ws.Range("H113").Value = 699.4
$ws.Range("J113").Value = 674.75
$ws.Range("L113").Value = 2024.25
$ws.Range("N113").Value = -6364.25
# Row 132
$ws.Range("H132").Value = 1158.25
$ws.Range("I132").Value = 1158.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10424.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7894.25
$ws.Range("N132").ClearContents()
# Row 139
$ws.Range("H139").Value = 3773.2083
$ws.Range("I139").Value = 1921.75
$ws.Range("J139").Value = 5624.6665
$ws.Range("K139").Value = 5765.25
$ws.Range("L139").Value = 16873.9995
$ws.Range("M139").Value = -625.25
$ws.Range("N139").Value = -27153.9995

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3724.8965
$ws.Range("J126").Value = 4587
$ws.Range("L126").Value = 13761
$ws.Range("N126").Value = -18701
# Row 132
$ws.Range("H132").Value = 10802
$ws.Range("I132").Value = 11002.5
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 33007.5
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -30477.5
$ws.Range("N132").Value = -35060

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5710.4443
$ws.Range("I7").Value = 4753.8887
$ws.Range("J7").Value = 6667
$ws.Range("K7").Value = 4753.8887
$ws.Range("L7").Value = 6667
$ws.Range("M7").Value = -4641.8887
$ws.Range("N7").Value = -6891
# Row 61
$ws.Range("H61").Value = 3665.8333
$ws.Range("I61").Value = 3399.8
$ws.Range("K61").Value = 3399.8
$ws.Range("M61").Value = -3197.8
# Row 113
$ws.Range("H113").Value = 3665.8333
$ws.Range("I113").Value = 3399.8
$ws.Range("K113").Value = 3399.8
$ws.Range("M113").Value = -1229.8
# Row 126
$ws.Range("H126").Value = 5710.4443
$ws.Range("I126").Value = 4753.8887
$ws.Range("J126").Value = 6667
$ws.Range("K126").Value = 14261.6661
$ws.Range("L126").Value = 20001
$ws.Range("M126").Value = -11791.6661
$ws.Range("N126").Value = -24941
# Row 132
$ws.Range("H132").Value = 4620.0464
$ws.Range("I132").Value = 3459.8125
$ws.Range("J132").Value = 7995.273
$ws.Range("K132").Value = 10379.4375
$ws.Range("L132").Value = 23985.819
$ws.Range("M132").Value = -7849.4375
$ws.Range("N132").Value = -29045.819

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 628.0769
$ws.Range("I113").Value = 695.1111
$ws.Range("K113").Value = 2085.3333
$ws.Range("M113").Value = 84.66670000000022
# Row 125
$ws.Range("H125").Value = 56666.668
$ws.Range("J125").Value = 56666.668
$ws.Range("L125").Value = 56666.668
$ws.Range("N125").Value = -66506.66800000001
# Row 132
$ws.Range("H132").Value = 50514790
$ws.Range("I132").Value = 18520814
$ws.Range("J132").Value = 62512532
$ws.Range("K132").Value = 55562442
$ws.Range("L132").Value = 187537596
$ws.Range("M132").Value = -55559912
$ws.Range("N132").Value = -187542656

